$wb = $excel.ActiveWorkbook

# Rename the "Include from ..." sheets to "Include #N"
$wb.Worksheets.Item("Include from MedComCorePracti").Name = "Include #0"
$wb.Worksheets.Item("Include from MedComCorePracti 2").Name = "Include #1"
$wb.Worksheets.Item("Include from NullFlavor").Name = "Include #2"

# Update the IG version number on the Metadata sheet (row with "Version" in col A)
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B3").Value = "1.7.1"
